# Auto-generated edit script: updates Leve-profit market-data columns (H-N)
# across all 8 job sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 7573
$ws.Range("I74").Value = 7463.231
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 7463.231
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -6527.231
$ws.Range("N74").Value = -10872
# Row 77
$ws.Range("H77").Value = 7573
$ws.Range("I77").Value = 7463.231
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 37316.155
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -32636.155
$ws.Range("N77").Value = -54360
# Row 86
$ws.Range("H86").Value = 6888.8335
$ws.Range("I86").Value = 4444.3335
$ws.Range("K86").Value = 4444.3335
$ws.Range("M86").Value = -3321.3335
# Row 89
$ws.Range("H89").Value = 6888.8335
$ws.Range("I89").Value = 4444.3335
$ws.Range("K89").Value = 22221.6675
$ws.Range("M89").Value = -16605.6675
# Row 98
$ws.Range("H98").Value = 8000
$ws.Range("I98").Value = 2000
$ws.Range("K98").Value = 2000
$ws.Range("M98").Value = -502
# Row 122
$ws.Range("H122").Value = 8000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
# Row 125
$ws.Range("H125").Value = 471
$ws.Range("I125").Value = 471
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 4239
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -1779
$ws.Range("N125").Value = ""
# Row 137
$ws.Range("H137").Value = 1702
$ws.Range("I137").Value = 1702
$ws.Range("K137").Value = 5106
$ws.Range("M137").Value = -2556

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 700
$ws.Range("I8").Value = 700
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 700
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -556
$ws.Range("N8").Value = ""
# Row 31
$ws.Range("H31").Value = 23333.334
$ws.Range("I31").Value = 10000
$ws.Range("K31").Value = 10000
$ws.Range("M31").Value = -9706
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""
# Row 45
$ws.Range("H45").Value = 4079.8
$ws.Range("I45").Value = 4599.75
$ws.Range("K45").Value = 4599.75
$ws.Range("M45").Value = -4222.75
# Row 61
$ws.Range("H61").Value = 1212.2727
$ws.Range("I61").Value = 1270.8889
$ws.Range("J61").Value = 948.5
$ws.Range("K61").Value = 1270.8889
$ws.Range("L61").Value = 948.5
$ws.Range("M61").Value = -1058.8889
$ws.Range("N61").Value = -1372.5
# Row 80
$ws.Range("H80").Value = 46997.5
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 73995
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 73995
$ws.Range("M80").Value = -19002
$ws.Range("N80").Value = -75991
# Row 83
$ws.Range("H83").Value = 46997.5
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 73995
$ws.Range("K83").Value = 60000
$ws.Range("L83").Value = 221985
$ws.Range("M83").Value = -55008
$ws.Range("N83").Value = -231969
# Row 122
$ws.Range("H122").Value = 3037.625
$ws.Range("I122").Value = 3214.2856
$ws.Range("J122").Value = 1801
$ws.Range("K122").Value = 9642.856800000001
$ws.Range("L122").Value = 5403
$ws.Range("M122").Value = -7192.856800000001
$ws.Range("N122").Value = -10303
# Row 136
$ws.Range("H136").Value = 1212.2727
$ws.Range("I136").Value = 1270.8889
$ws.Range("J136").Value = 948.5
$ws.Range("K136").Value = 3812.6667
$ws.Range("L136").Value = 2845.5
$ws.Range("M136").Value = -1262.6667
$ws.Range("N136").Value = -7945.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 141
$ws.Range("H141").Value = 179997
$ws.Range("J141").Value = 199998
$ws.Range("L141").Value = 199998
$ws.Range("N141").Value = -210358

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""
# Row 58
$ws.Range("H58").Value = 6531.6665
$ws.Range("I58").Value = 2297.5
$ws.Range("K58").Value = 2297.5
$ws.Range("M58").Value = -2094.5
# Row 105
$ws.Range("H105").Value = 5598.2
$ws.Range("I105").Value = 7250
$ws.Range("J105").Value = 4497
$ws.Range("K105").Value = 7250
$ws.Range("L105").Value = 4497
$ws.Range("M105").Value = -5503
$ws.Range("N105").Value = -7991
# Row 122
$ws.Range("H122").Value = 5149
$ws.Range("I122").Value = 7165
$ws.Range("K122").Value = 21495
$ws.Range("M122").Value = -19045
# Row 134
$ws.Range("H134").Value = 11427.571
$ws.Range("I134").Value = 11663.167
$ws.Range("K134").Value = 34989.501
$ws.Range("M134").Value = -32454.501
# Row 136
$ws.Range("H136").Value = 6531.6665
$ws.Range("I136").Value = 2297.5
$ws.Range("K136").Value = 6892.5
$ws.Range("M136").Value = -4342.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 14064
$ws.Range("I4").Value = 17430
$ws.Range("K4").Value = 52290
$ws.Range("M4").Value = -52178
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""
# Row 22
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 2997
$ws.Range("M22").Value = -2828
# Row 27
$ws.Range("H27").Value = 999
$ws.Range("I27").Value = 999
$ws.Range("K27").Value = 2997
$ws.Range("M27").Value = -2895
# Row 34
$ws.Range("H34").Value = 5830.5713
$ws.Range("J34").Value = 7600
$ws.Range("L34").Value = 22800
$ws.Range("N34").Value = -22968
# Row 113
$ws.Range("H113").Value = 1211.5333
$ws.Range("I113").Value = 722.5
$ws.Range("J113").Value = 1389.3636
$ws.Range("K113").Value = 2167.5
$ws.Range("L113").Value = 4168.0908
$ws.Range("M113").Value = 2.5
$ws.Range("N113").Value = -8508.0908

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""
# Row 70
$ws.Range("H70").Value = 12000
$ws.Range("J70").Value = 12000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12540
# Row 73
$ws.Range("H73").Value = 12000
$ws.Range("J73").Value = 12000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -13872
# Row 80
$ws.Range("H80").Value = 1975
$ws.Range("I80").Value = 1950
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1950
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -952
$ws.Range("N80").Value = -3996
# Row 83
$ws.Range("H83").Value = 1975
$ws.Range("I83").Value = 1950
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 9750
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -4758
$ws.Range("N83").Value = -19984
# Row 95
$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 100000
$ws.Range("N95").Value = -105492

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2566
$ws.Range("I68").Value = 2679.2
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2679.2
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1930.2
$ws.Range("N68").Value = -3498
# Row 71
$ws.Range("H71").Value = 2566
$ws.Range("I71").Value = 2679.2
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 13396
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -9652
$ws.Range("N71").Value = -17488
# Row 130
$ws.Range("H130").Value = 58333
$ws.Range("J130").Value = 58333
$ws.Range("L130").Value = 58333
$ws.Range("N130").Value = -68373
# Row 132
$ws.Range("H132").Value = 5663.8335
$ws.Range("I132").Value = 4995.5
$ws.Range("K132").Value = 14986.5
$ws.Range("M132").Value = -12456.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 6497.25
$ws.Range("I62").Value = 3300
$ws.Range("K62").Value = 3300
$ws.Range("M62").Value = -2676
# Row 65
$ws.Range("H65").Value = 6497.25
$ws.Range("I65").Value = 3300
$ws.Range("K65").Value = 16500
$ws.Range("M65").Value = -13380
# Row 96
$ws.Range("H96").Value = 6198
$ws.Range("I96").Value = 6198
$ws.Range("K96").Value = 6198
$ws.Range("M96").Value = -4825
# Row 132
$ws.Range("H132").Value = 2103.85
$ws.Range("I132").Value = 1139
$ws.Range("K132").Value = 3417
$ws.Range("M132").Value = -887
# Row 136
$ws.Range("H136").Value = 6712.1665
$ws.Range("I136").Value = 5727.4443
$ws.Range("K136").Value = 17182.3329
$ws.Range("M136").Value = -14632.3329
